$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# --- Add the new Input/Output directory settings rows to the Settings sheet ---
$wsSettings.Range("A6").Value = "OutputDirectory"
$wsSettings.Range("B6").Value = "Data\Output\SampleOutput.csv"
$wsSettings.Range("C6").Value = "Output file path"

$wsSettings.Range("A7").Value = "InputDirectory"
$wsSettings.Range("B7").Value = "Data\Input\SampleInput.csv"
$wsSettings.Range("C7").Value = "Input file path"

# Row height tweak that comes from re-wrapping the description text (row 4)
$wsSettings.Rows.Item(4).RowHeight = 28.8
$wsConstants.Rows.Item(2).RowHeight = 28.8

# --- Update selection / active sheet state ---
# Assets sheet no longer the active tab; move its selection to A2
$wsAssets.Activate() | Out-Null
$wsAssets.Range("A2").Select() | Out-Null

# Constants sheet selection left untouched

# Settings becomes the active tab, selection moves down to C8
$wsSettings.Activate() | Out-Null
$wsSettings.Range("C8").Select() | Out-Null
